$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) values are textual (often multi-dot like "26.867.06")
# in the source sheet, not real numbers. Force text storage so Excel
# COM does not auto-coerce single-dot-looking values (e.g. "309.82")
# into floating point numbers, then restore the default style.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.867.06"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.84%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.809.39"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.09%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4640"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.86%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3703"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.69%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07354"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.82%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8769"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.17%  "
$ws.Range("E11").Value = "  -1.98%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.768.92"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.42%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.357"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.36%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.516"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.93%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.83"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.31%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.07048"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.32%  "
$ws.Range("E17").Value = "  +0.12%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008689"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.56%  "
$ws.Range("E19").Value = "  +0.11%  "
$ws.Range("E20").Value = "  -2.42%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "26.889.35"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.77%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.313"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.70%  "
$ws.Range("E23").Value = "  -3.36%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.051.07"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.20%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.901"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.97%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.62"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.26%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.39"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.13%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.160"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.35%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.325"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.53%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "116.05"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.05%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08910"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.04%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7531"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.26%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.158"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.36%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.922"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.37%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.459"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.15%  "
$ws.Range("E36").Value = "  +0.10%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.101"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.35%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01969"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.59%  "
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05252"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.52%  "
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.430"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.83%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.929"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.87%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5324"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.45%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.178"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.47%  "
$ws.Range("E44").Value = "  -2.34%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.489"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.37%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4980"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.26%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.30"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.55%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "103.92"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.39%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.669"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.24%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06297"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.44%  "
